$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: columns B and C swap meaning (B1 = Ano, C1 = Variável)
$ws.Cells.Item(1, 2).Value = "Ano"
$ws.Cells.Item(1, 3).Value = "Variável"

# Data rows 2-19: reorganized so Brasil comes first, then Nordeste, then Sergipe.
# Column B now holds the date (dd/mm/yyyy as literal text) and column C holds the category.
# The leading apostrophe forces the date-looking text to stay text instead of
# being auto-converted into a date serial by Excel.
$data = @(
    @("Brasil",   "'01/01/2019", "Feminicídio", 1.491603669709312, $null),
    @("Brasil",   "'01/01/2020", "Feminicídio", 1.53032839946819,  $null),
    @("Brasil",   "'01/01/2021", "Feminicídio", 1.569654974814453, $null),
    @("Brasil",   "'01/01/2022", "Feminicídio", 1.532101471544391, $null),
    @("Brasil",   "'01/01/2023", "Feminicídio", 1.503723149276654, $null),
    @("Brasil",   "'01/01/2024", "Feminicídio", 1.443586697474013, $null),
    @("Nordeste", "'01/01/2019", "Feminicídio", 1.497286779739304, $null),
    @("Nordeste", "'01/01/2020", "Feminicídio", 1.421470954921448, $null),
    @("Nordeste", "'01/01/2021", "Feminicídio", 1.450454625600147, $null),
    @("Nordeste", "'01/01/2022", "Feminicídio", 1.326824150475039, $null),
    @("Nordeste", "'01/01/2023", "Feminicídio", 1.33400727814508,  $null),
    @("Nordeste", "'01/01/2024", "Feminicídio", 1.31935864980953,  $null),
    @("Sergipe",  "'01/01/2019", "Feminicídio", 1.766753237575308, 7),
    @("Sergipe",  "'01/01/2020", "Feminicídio", 1.166870869068754, 20),
    @("Sergipe",  "'01/01/2021", "Feminicídio", 1.651913700724447, 9),
    @("Sergipe",  "'01/01/2022", "Feminicídio", 1.555611048113412, 12),
    @("Sergipe",  "'01/01/2023", "Feminicídio", 1.298940632730219, 19),
    @("Sergipe",  "'01/01/2024", "Feminicídio", 0.8052392083854389, 25)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}
